$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header label updates (switch RQ134 from DRE-based to distance-based naming) ---
$ws.Range("O1").Value = "P1E/Population Alignment score (ADM|target)"
$ws.Range("P1").Value = "DRE/Distance Alignment score (ADM|target)"
$ws.Range("Q1").Value = "P1E/Population Alignment score (Delegator|target)"
$ws.Range("R1").Value = "DRE/Distance Alignment score (Delegator|target)"

$ws.Range("Y1").Value = "P1E/Population Alignment score (Delegator|Observed_ADM (target))"
$ws.Range("Z1").Value = "DRE/Distance Alignment score (Delegator|Observed_ADM (target))"

# --- Row 3 definition updates ---
$ws.Range("O3").Value = "Calculated alignment score between the KDMA measurement of the ADM aligned to a particular target and that target using the population-based endpoints"
$ws.Range("P3").Value = "Calculated alignment score between the KDMA measurement of the ADM aligned to a particular target and that target using the distance-based endpoints in the Phase 1 server for ADEPT, no change for ST"
$ws.Range("Q3").Value = "Calculated alignment score between the KDMA measurement of a delegator and a target using the population-based endpoints"
$ws.Range("R3").Value = "Calculated alignment score between the KDMA measurement of a delegator and a target using the distance-based endpoints in the Phase 1 server for ADEPT, no change for ST"

$ws.Range("Y3").Value = "Compares the KDMA measurement based on delegator responses to the subset of probes from an ADM they observed in the delegation survey using the Phase 1 server and new population-based endpoints"
$ws.Range("Z3").Value = "Compares the KDMA measurement based on delegator responses to the subset of probes from an ADM they observed in the delegation survey using the distance-based endpoints in the Phase 1 server for ADEPT, no change for ST"

$wb.Save()
